$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = "codeforiati:group-name", Column E = "codeforiati:group-code" (before edit).
# The edit swaps these two columns so D becomes the code and E becomes the name,
# matching the reordering of the shared-string table in the source diff.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $eVal
    $ws.Cells.Item($r, 5).Value2 = $dVal
}
